$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text format to preserve exact string
# representation (values like "220.32" would otherwise be auto-converted
# to numbers by Excel's type inference).
$dCells = @('D2', 'D3', 'D5', 'D6', 'D8', 'D9', 'D10', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D25', 'D26', 'D27', 'D30', 'D31', 'D33', 'D34', 'D35', 'D37', 'D39', 'D41', 'D45', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($c in $dCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range('D2').Value = '33.890.42'
$ws.Range('D3').Value = '1.768.42'
$ws.Range('D5').Value = '220.32'
$ws.Range('D6').Value = '0.545'
$ws.Range('D8').Value = '30.84'
$ws.Range('D9').Value = '0.283'
$ws.Range('D10').Value = '0.0703'
$ws.Range('D12').Value = '2.026.77'
$ws.Range('D13').Value = '1.765.34'
$ws.Range('D14').Value = '10.46'
$ws.Range('D15').Value = '0.621'
$ws.Range('D16').Value = '33.915.01'
$ws.Range('D17').Value = '4.18'
$ws.Range('D18').Value = '67.49'
$ws.Range('D19').Value = '242.70'
$ws.Range('D20').Value = '0.0₃0772'
$ws.Range('D21').Value = '1.00'
$ws.Range('D22').Value = '10.49'
$ws.Range('D25').Value = '157.21'
$ws.Range('D26').Value = '16.29'
$ws.Range('D27').Value = '6.95'
$ws.Range('D30').Value = '0.0519'
$ws.Range('D31').Value = '3.69'
$ws.Range('D33').Value = '3.49'
$ws.Range('D34').Value = '1.80'
$ws.Range('D35').Value = '1.394.76'
$ws.Range('D37').Value = '0.630'
$ws.Range('D39').Value = '0.924'
$ws.Range('D41').Value = '78.33'
$ws.Range('D45').Value = '0.0487'
$ws.Range('D46').Value = '1.04'
$ws.Range('D47').Value = '1.922.12'
$ws.Range('D48').Value = '103.78'
$ws.Range('D50').Value = '11.78'
$ws.Range('D51').Value = '0.0₆0120'

foreach ($c in $dCells) { $ws.Range($c).Style = "Normal" }

# Column E (Volume 1h) updates - plain text values (already safe; contain
# "%" and padding spaces so Excel keeps them as text).
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  -5.73%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  +5.15%  '
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('E14').Value = '  -5.54%  '
$ws.Range('E15').Value = '  -2.19%  '
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('E17').Value = '  -2.32%  '
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('E19').Value = '  -4.51%  '
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('E23').Value = '  -4.73%  '
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('E28').Value = '  -2.55%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('E35').Value = '  -3.69%  '
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('E39').Value = '  +2.85%  '
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('E41').Value = '  -5.60%  '
$ws.Range('E42').Value = '  -5.13%  '
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('E45').Value = '  -3.80%  '
$ws.Range('E46').Value = '  -1.45%  '
$ws.Range('E47').Value = '  -2.00%  '
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('E51').Value = '  -1.72%  '
